{"js": "// Replace the 25 \"AA\u00d7BB=\" multiplication prompts in the table with their\n// new values, in document order. One source string (\"23\u00d793=\") occurs twice\n// with two different replacement targets, so occurrences are consumed in\n// document order rather than via a single blind global replace.\n\nconst pairs = [\n  [\"60\u00d717=\", \"46\u00d715=\"],\n  [\"89\u00d794=\", \"65\u00d788=\"],\n  [\"14\u00d776=\", \"36\u00d750=\"],\n  [\"90\u00d788=\", \"17\u00d743=\"],\n  [\"59\u00d721=\", \"24\u00d724=\"],\n  [\"60\u00d726=\", \"89\u00d759=\"],\n  [\"23\u00d793=\", \"53\u00d736=\"],\n  [\"87\u00d775=\", \"73\u00d720=\"],\n  [\"66\u00d732=\", \"24\u00d792=\"],\n  [\"70\u00d797=\", \"50\u00d728=\"],\n  [\"52\u00d782=\", \"45\u00d753=\"],\n  [\"28\u00d755=\", \"90\u00d793=\"],\n  [\"47\u00d770=\", \"97\u00d799=\"],\n  [\"62\u00d766=\", \"23\u00d745=\"],\n  [\"90\u00d712=\", \"87\u00d718=\"],\n  [\"59\u00d779=\", \"89\u00d767=\"],\n  [\"23\u00d793=\", \"70\u00d761=\"],\n  [\"98\u00d737=\", \"59\u00d755=\"],\n  [\"77\u00d781=\", \"20\u00d782=\"],\n  [\"47\u00d787=\", \"68\u00d787=\"],\n  [\"83\u00d733=\", \"21\u00d784=\"],\n  [\"45\u00d794=\", \"39\u00d764=\"],\n  [\"88\u00d734=\", \"29\u00d712=\"],\n  [\"11\u00d711=\", \"73\u00d783=\"],\n  [\"60\u00d711=\", \"27\u00d719=\"],\n];\n\n// Each replacement is applied immediately, so once the first occurrence of\n// a repeated source string (e.g. \"23\u00d793=\") is rewritten, a fresh search for\n// that same source string only finds the still-unprocessed occurrence(s) \u2014\n// meaning we always take result index 0, processed strictly in document\n// (pairs-array) order.\nfor (const [from, to] of pairs) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find \"${from}\" to replace with \"${to}\"`);\n  }\n\n  results.items[0].insertText(to, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"AA\u00d7BB=\" multiplication prompts in the table with their\n# new values, in document order. One source string (\"23\u00d793=\") occurs twice\n# with two different replacement targets, so each replacement re-searches\n# the whole document and replaces only the first (wdReplaceOne) remaining\n# match, which \u2014 since earlier matches are already rewritten \u2014 always lands\n# on the correct, next-in-document-order occurrence.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"60\u00d717=\", \"46\u00d715=\"),\n    @(\"89\u00d794=\", \"65\u00d788=\"),\n    @(\"14\u00d776=\", \"36\u00d750=\"),\n    @(\"90\u00d788=\", \"17\u00d743=\"),\n    @(\"59\u00d721=\", \"24\u00d724=\"),\n    @(\"60\u00d726=\", \"89\u00d759=\"),\n    @(\"23\u00d793=\", \"53\u00d736=\"),\n    @(\"87\u00d775=\", \"73\u00d720=\"),\n    @(\"66\u00d732=\", \"24\u00d792=\"),\n    @(\"70\u00d797=\", \"50\u00d728=\"),\n    @(\"52\u00d782=\", \"45\u00d753=\"),\n    @(\"28\u00d755=\", \"90\u00d793=\"),\n    @(\"47\u00d770=\", \"97\u00d799=\"),\n    @(\"62\u00d766=\", \"23\u00d745=\"),\n    @(\"90\u00d712=\", \"87\u00d718=\"),\n    @(\"59\u00d779=\", \"89\u00d767=\"),\n    @(\"23\u00d793=\", \"70\u00d761=\"),\n    @(\"98\u00d737=\", \"59\u00d755=\"),\n    @(\"77\u00d781=\", \"20\u00d782=\"),\n    @(\"47\u00d787=\", \"68\u00d787=\"),\n    @(\"83\u00d733=\", \"21\u00d784=\"),\n    @(\"45\u00d794=\", \"39\u00d764=\"),\n    @(\"88\u00d734=\", \"29\u00d712=\"),\n    @(\"11\u00d711=\", \"73\u00d783=\"),\n    @(\"60\u00d711=\", \"27\u00d719=\")\n)\n\nforeach ($pair in $pairs) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $result = $find.Execute($from, $false, $false, $false, $false, $false, $true, 1, $false, $to, 1)\n    if (-not $result) {\n        Write-Output \"WARNING: could not find '$from' to replace with '$to'\"\n    }\n}\n"}
